$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$dateFmt = $ws.Range("B14").NumberFormat

$ws.Range("A15").Value = "estudio db40 pero lo desecho finalmente, ya no existe. Inicio el programaa con la ventana de login en lugar de la principal"
$ws.Range("B15").NumberFormat = $dateFmt
$ws.Range("B15").Value = (Get-Date -Year 2020 -Month 8 -Day 12).Date

$ws.Range("A16").Value = "estudio ObjectDB - tutorial"
$ws.Range("B16").NumberFormat = $dateFmt
$ws.Range("B16").Value = (Get-Date -Year 2020 -Month 8 -Day 13).Date

$ws.Range("B17").Select()
